$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these price cells remain text (not auto-converted to numbers)
$textCells = @("D5","D10","D11","D16","D19","D21","D22","D24","D25","D30","D32","D36","D42","D43","D44","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '26.666.41'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').Value = '1.594.18'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '211.46'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').Value = '19.67'
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('D11').Value = '0.0836'
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').Value = '1.818.01'
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('D13').Value = '1.574.82'
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('E15').Value = '  -3.29%  '
$ws.Range('D16').Value = '64.71'
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('D17').Value = '26.640.00'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').Value = '209.13'
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -2.39%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D24').Value = '8.86'
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').Value = '146.66'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  -1.28%  '
$ws.Range('D30').Value = '0.0503'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('D32').Value = '3.21'
$ws.Range('E32').Value = '  -3.85%  '
$ws.Range('E33').Value = '  -8.48%  '
$ws.Range('E34').Value = '  -2.61%  '
$ws.Range('D35').Value = '1.284.84'
$ws.Range('E35').Value = '  -5.69%  '
$ws.Range('D36').Value = '2.43'
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('E37').Value = '  -5.66%  '
$ws.Range('E38').Value = '  -2.94%  '
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '2.20'
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.36'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').Value = '63.52'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('D45').Value = '1.730.35'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('D46').Value = '89.69'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '0.873'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.63'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').Value = '0.0981'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('E50').Value = '  -1.90%  '
$ws.Range('D51').Value = '7.51'
$ws.Range('E51').Value = '  -1.50%  '
